$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 374
$ws1.Range("F3").Value = 72
$ws1.Range("F4").Value = 292
$ws1.Range("F5").Value = 4273
$ws1.Range("F6").Value = 41
$ws1.Range("F7").Value = 459

# Sheet "全部类型" (All types) - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 374
$ws4.Range("F3").Value = 72
$ws4.Range("F4").Value = 292
$ws4.Range("F5").Value = 4273
$ws4.Range("F8").Value = 41
$ws4.Range("F9").Value = 459
